# "update tasks for Tet Holiday"
# Three tasks (rows 22, 27, 28 on Sheet1) move from "waiting" to
# "on processing" and get an owner ("Hung") assigned in column C.
# The view is also scrolled down a bit with a new active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Assign owner "Hung" to the three newly-active tasks (column C currently blank)
$ws.Range("C22").Value = "Hung"
$ws.Range("C27").Value = "Hung"
$ws.Range("C28").Value = "Hung"

# Move status from "waiting" to "on processing" for the same rows, copying the
# cell format (fill colour etc.) used by the existing "on processing" cells
# (column G, rows 23/24) so the new cells render identically to their peers.
$ws.Range("G23").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G28").PasteSpecial(-4122)

$ws.Range("G22").Value = "on processing"
$ws.Range("G27").Value = "on processing"
$ws.Range("G28").Value = "on processing"

# Update the view: scroll further down the list and move the active selection.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C35").Select() | Out-Null
